# Update cryptos list (values scraped from coinranking.com)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.921.39'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '1.670.39'
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.85'
$ws.Range("E5").Value = '  -0.13%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.516'
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0621'
$ws.Range("E8").Value = '  +0.81%  '
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.250'
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("E11").Value = '  +1.70%  '
$ws.Range("D12").Value = '1.905.31'
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").Value = '1.696.97'
$ws.Range("E13").Value = '  +2.48%  '
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("E15").Value = '  +1.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.46'
$ws.Range("E16").Value = '  +0.49%  '
$ws.Range("D17").Value = '26.913.35'
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.02'
$ws.Range("E18").Value = '  +3.85%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '233.10'
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.42'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.15'
$ws.Range("E23").Value = '  -1.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").Value = '  -2.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.19'
$ws.Range("E25").Value = '  +0.54%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.94'
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("E29").Value = '  -1.96%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("D33").Value = '1.457.09'
$ws.Range("E33").Value = '  -5.75%  '
$ws.Range("E34").Value = '  +1.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.63'
$ws.Range("E35").Value = '  +1.79%  '
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("E37").Value = '  -0.16%  '
$ws.Range("E38").Value = '  +0.72%  '
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("E40").Value = '  +13.53%  '
$ws.Range("E41").Value = '  -4.45%  '
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("E43").Value = '  +2.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.25'
$ws.Range("E44").Value = '  +0.75%  '
$ws.Range("D45").Value = '1.811.96'
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("E46").Value = '  +0.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.71'
$ws.Range("E47").Value = '  +0.40%  '
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0103'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.102'
$ws.Range("E50").Value = '  +2.80%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0507'
$ws.Range("E51").Value = '  +0.29%  '
